$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,9,11,13,14,15,16,17)
$cols = @("B","C","D")

$c1 = [char]0xC2
$c2 = [char]0xB1
$bad = "$c1$c2"
$good = [string][char]0xB1

foreach ($r in $rows) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $val = $cell.Value()
        if ($val -ne $null) {
            $newVal = $val.Replace($bad, $good)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
